$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: new A22 cell
$ws.Range("A22").Formula = "=A21+B21+D21"

# Row 25: new A25, B25 (register shared string "sec" first)
$ws.Range("A25").Formula = "=A21/1000"
$ws.Range("B25").Value = "sec"

# Row 26: new A26, B26 (register shared string "min" second)
$ws.Range("A26").Formula = "=TRUNC(A25/60)"
$ws.Range("B26").Value = "min"

# Row 27: new A27, B27 (reuse shared string "sec")
$ws.Range("A27").Formula = "=MOD(A25,60)"
$ws.Range("A27").NumberFormat = "0"
$ws.Range("B27").Value = "sec"

# Row 23: new C23 (label, registers shared string "Throughput" last) and updated D23 formula
$ws.Range("C23").Value = "Throughput"
$ws.Range("C23").HorizontalAlignment = -4152  # xlRight
$ws.Range("D23").Formula = "=1-D22/A22"

# Update selection to D23
$ws.Range("D23").Select()
